$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 244; this shifts the existing rows 244-355
# down to 245-356 (the old last row, 355, duplicates into the new 356,
# matching the target diff), and the worksheet's used range / dimension
# grows to A1:R356 automatically.
$ws.Rows(244).Insert()

# Populate the newly inserted row 244 with the new record.
$ws.Range("A244").Value = 8
$ws.Range("B244").Value = "Terminal La Palmera de La Serena"
$ws.Range("C244").Value = "Coquimbo"
$ws.Range("D244").Value = 44452
$ws.Range("E244").Value = 4
$ws.Range("F244").Value = 100112024
$ws.Range("G244").Value = "Choclo"
$ws.Range("H244").Value = "Dulce o Americano"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 500
$ws.Range("K244").Value = 34000
$ws.Range("L244").Value = 35000
$ws.Range("M244").Value = 34500
$ws.Range("N244").Value = "$/malla 70 unidades"
$ws.Range("O244").Value = "Región de Arica y Parinacota"
$ws.Range("P244").Value = 493
$ws.Range("Q244").Value = 70
$ws.Range("R244").Value = "Hortaliza"
